$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (pushes the old rows 4-8 down to 5-9,
# carrying their existing formatting/values along with them).
$ws.Rows.Item(4).Insert()

# Populate the freshly inserted row 4 with the new weekly price record.
$ws.Range("A4").Value = 8
$ws.Range("B4").Value = "Terminal La Palmera de La Serena"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44498
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 100112013
$ws.Range("G4").Value = "Alcachofa"
$ws.Range("H4").Value = "Española"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 8500
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8750
$ws.Range("N4").Value = "$/caja 30 unidades"
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 292
$ws.Range("Q4").Value = 30
$ws.Range("R4").Value = "Hortaliza"
